$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Slicer Used" column (column E) entirely; remaining columns shift left
$ws.Range("E1").EntireColumn.Delete()

# Fill in "Printed By" (column C) for existing rows 2-9
$ws.Range("C2").Value = "Nathan"
$ws.Range("C3").Value = "Nathan"
$ws.Range("C4").Value = "Nathan"
$ws.Range("C5").Value = "Nathan"
$ws.Range("C6").Value = "Nathan"
$ws.Range("C7").Value = "Nathan"
$ws.Range("C8").Value = "Nathan"
$ws.Range("C9").Value = "Nathan"

# Fill in "Part Fit Tested?" (column F) for existing rows 2-9
$ws.Range("F2").Value = "sort of"
$ws.Range("F3").Value = "yes"
$ws.Range("F4").Value = "yes"
$ws.Range("F5").Value = "yes"
$ws.Range("F6").Value = "yes"
$ws.Range("F7").Value = "yes"
$ws.Range("F8").Value = "yes"
$ws.Range("F9").Value = "yes"

# Add two new parts as additional rows
$ws.Range("A10").Value = 6
$ws.Range("B10").Value = "main wrist new"

$ws.Range("A11").Value = 6
$ws.Range("B11").Value = "tricep"

# Update selection to reflect where the user left off
$ws.Range("C10").Select()
